# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, the
# following "(c) 2020 . Contact: ..." paragraph, and the blank paragraph
# that precedes them (right after the page-break paragraph that follows
# "Requisitos"). Paragraph indices are 1-based in the Word object model.
#
# Before (1-indexed paragraphs near the end of the document):
#   54: "Requisitos"
#   55: (empty)
#   56: (empty, pageBreakBefore)
#   57: (empty)                                            <- delete
#   58: "Ver no Jupiter Salvar em pdf Salvar em docx"       <- delete
#   59: "(c) 2020 . Contact: ..."                           <- delete
#   60: (empty)
#   61: (empty, pageBreakBefore)

$d = $word.ActiveDocument

# Delete from the bottom up so earlier indices stay valid.
$d.Paragraphs.Item(59).Range.Delete()
$d.Paragraphs.Item(58).Range.Delete()
$d.Paragraphs.Item(57).Range.Delete()
